# Renamed few transcripts. Updated the DataSheet
# Replace speaker names in column D: "Davis" -> "T", "Student" -> "S"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)  # Column D
    $val = $cell.Text
    if ($val -eq "Davis") {
        $cell.Value = "T"
    } elseif ($val -eq "Student") {
        $cell.Value = "S"
    }
}
